# Updates currentAveragePrice / Leve price / profit columns (H-N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets,
# reflecting refreshed market-board data from the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 201
$ws.Cells.Item(6, 9).Value = 201
$ws.Cells.Item(6, 11).Value = 603
$ws.Cells.Item(6, 13).Value = -491
$ws.Cells.Item(7, 8).Value = 1700
$ws.Cells.Item(7, 9).Value = 1700
$ws.Cells.Item(7, 11).Value = 1700
$ws.Cells.Item(7, 13).Value = -1588
$ws.Cells.Item(8, 8).Value = 1360
$ws.Cells.Item(8, 9).Value = 206.66667
$ws.Cells.Item(8, 10).Value = 3666.6667
$ws.Cells.Item(8, 11).Value = 620.00001
$ws.Cells.Item(8, 12).Value = 11000.0001
$ws.Cells.Item(8, 13).Value = -481.00001
$ws.Cells.Item(8, 14).Value = -11278.0001
$ws.Cells.Item(10, 8).Value = 4600
$ws.Cells.Item(10, 10).Value = 5000
$ws.Cells.Item(10, 12).Value = 5000
$ws.Cells.Item(10, 14).Value = -5586
$ws.Cells.Item(14, 8).Value = 1700
$ws.Cells.Item(14, 9).Value = 1700
$ws.Cells.Item(14, 11).Value = 1700
$ws.Cells.Item(14, 13).Value = -1509
$ws.Cells.Item(28, 8).Value = 1030.7
$ws.Cells.Item(28, 9).Value = 914.26666
$ws.Cells.Item(28, 10).Value = 1380
$ws.Cells.Item(28, 11).Value = 914.26666
$ws.Cells.Item(28, 12).Value = 1380
$ws.Cells.Item(28, 13).Value = -429.26666
$ws.Cells.Item(28, 14).Value = -2350
$ws.Cells.Item(31, 8).Value = 500
$ws.Cells.Item(31, 9).Value = 500
$ws.Cells.Item(31, 10).Value = 0
$ws.Cells.Item(31, 11).Value = 1500
$ws.Cells.Item(31, 12).Value = 0
$ws.Cells.Item(31, 13).Value = -1270
$ws.Cells.Item(31, 14).Value = ""
$ws.Cells.Item(107, 8).Value = 1101.6538
$ws.Cells.Item(107, 9).Value = 1009.1905
$ws.Cells.Item(107, 10).Value = 1490
$ws.Cells.Item(107, 11).Value = 1009.1905
$ws.Cells.Item(107, 12).Value = 1490
$ws.Cells.Item(107, 13).Value = 910.8095
$ws.Cells.Item(107, 14).Value = -5330
$ws.Cells.Item(138, 8).Value = 4768.875
$ws.Cells.Item(138, 9).Value = 4705.8335
$ws.Cells.Item(138, 10).Value = 4806.7
$ws.Cells.Item(138, 11).Value = 14117.5005
$ws.Cells.Item(138, 12).Value = 14420.1
$ws.Cells.Item(138, 13).Value = -8977.500499999998
$ws.Cells.Item(138, 14).Value = -24700.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(12, 8).Value = 12500
$ws.Cells.Item(12, 10).Value = 12500
$ws.Cells.Item(12, 12).Value = 12500
$ws.Cells.Item(12, 14).Value = -12846
$ws.Cells.Item(45, 8).Value = 54227
$ws.Cells.Item(45, 9).Value = 68066.60000000001
$ws.Cells.Item(45, 10).Value = 2328.5
$ws.Cells.Item(45, 11).Value = 68066.60000000001
$ws.Cells.Item(45, 12).Value = 2328.5
$ws.Cells.Item(45, 13).Value = -67689.60000000001
$ws.Cells.Item(45, 14).Value = -3082.5
$ws.Cells.Item(115, 8).Value = 57663.125
$ws.Cells.Item(115, 9).Value = 0
$ws.Cells.Item(115, 10).Value = 57663.125
$ws.Cells.Item(115, 11).Value = 0
$ws.Cells.Item(115, 12).Value = 57663.125
$ws.Cells.Item(115, 13).Value = ""
$ws.Cells.Item(115, 14).Value = -60797.125
$ws.Cells.Item(122, 8).Value = 12196967
$ws.Cells.Item(122, 9).Value = 18520094
$ws.Cells.Item(122, 11).Value = 55560282
$ws.Cells.Item(122, 13).Value = -55557832
$ws.Cells.Item(124, 8).Value = 30000
$ws.Cells.Item(124, 10).Value = 30000
$ws.Cells.Item(124, 12).Value = 30000
$ws.Cells.Item(124, 14).Value = -39820
$ws.Cells.Item(125, 8).Value = 54979
$ws.Cells.Item(125, 10).Value = 54979
$ws.Cells.Item(125, 12).Value = 54979
$ws.Cells.Item(125, 14).Value = -64819

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 39140.344
$ws.Cells.Item(134, 9).Value = 43502.69
$ws.Cells.Item(134, 10).Value = 1333.3334
$ws.Cells.Item(134, 11).Value = 130508.07
$ws.Cells.Item(134, 12).Value = 4000.0002
$ws.Cells.Item(134, 13).Value = -127973.07
$ws.Cells.Item(134, 14).Value = -9070.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(23, 8).Value = 2947.3684
$ws.Cells.Item(27, 8).Value = 2947.3684
$ws.Cells.Item(33, 8).Value = 15490
$ws.Cells.Item(33, 9).Value = 5030
$ws.Cells.Item(33, 10).Value = 25950
$ws.Cells.Item(33, 11).Value = 5030
$ws.Cells.Item(33, 12).Value = 25950
$ws.Cells.Item(33, 13).Value = -4651
$ws.Cells.Item(33, 14).Value = -26708
$ws.Cells.Item(58, 8).Value = 3545.6453
$ws.Cells.Item(58, 9).Value = 2014.75
$ws.Cells.Item(58, 10).Value = 4512.5264
$ws.Cells.Item(58, 11).Value = 2014.75
$ws.Cells.Item(58, 12).Value = 4512.5264
$ws.Cells.Item(58, 13).Value = -1811.75
$ws.Cells.Item(58, 14).Value = -4918.5264
$ws.Cells.Item(62, 8).Value = 3733.3333
$ws.Cells.Item(62, 9).Value = 2500
$ws.Cells.Item(62, 10).Value = 4350
$ws.Cells.Item(62, 11).Value = 2500
$ws.Cells.Item(62, 12).Value = 4350
$ws.Cells.Item(62, 13).Value = -1876
$ws.Cells.Item(62, 14).Value = -5598
$ws.Cells.Item(65, 8).Value = 3733.3333
$ws.Cells.Item(65, 9).Value = 2500
$ws.Cells.Item(65, 10).Value = 4350
$ws.Cells.Item(65, 11).Value = 12500
$ws.Cells.Item(65, 12).Value = 21750
$ws.Cells.Item(65, 13).Value = -9380
$ws.Cells.Item(65, 14).Value = -27990
$ws.Cells.Item(132, 8).Value = 3043.516
$ws.Cells.Item(132, 9).Value = 3109.7334
$ws.Cells.Item(132, 10).Value = 2981.4375
$ws.Cells.Item(132, 11).Value = 9329.200199999999
$ws.Cells.Item(132, 12).Value = 8944.3125
$ws.Cells.Item(132, 13).Value = -6799.200199999999
$ws.Cells.Item(132, 14).Value = -14004.3125
$ws.Cells.Item(136, 8).Value = 3545.6453
$ws.Cells.Item(136, 9).Value = 2014.75
$ws.Cells.Item(136, 10).Value = 4512.5264
$ws.Cells.Item(136, 11).Value = 6044.25
$ws.Cells.Item(136, 12).Value = 13537.5792
$ws.Cells.Item(136, 13).Value = -3494.25
$ws.Cells.Item(136, 14).Value = -18637.5792

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 590.8077
$ws.Cells.Item(4, 9).Value = 266.3125
$ws.Cells.Item(4, 10).Value = 1110
$ws.Cells.Item(4, 11).Value = 798.9375
$ws.Cells.Item(4, 12).Value = 3330
$ws.Cells.Item(4, 13).Value = -686.9375
$ws.Cells.Item(4, 14).Value = -3554
$ws.Cells.Item(26, 8).Value = 125
$ws.Cells.Item(26, 9).Value = 100
$ws.Cells.Item(26, 10).Value = 150
$ws.Cells.Item(26, 11).Value = 300
$ws.Cells.Item(26, 12).Value = 450
$ws.Cells.Item(26, 13).Value = -12
$ws.Cells.Item(26, 14).Value = -1026
$ws.Cells.Item(113, 8).Value = 590.55554
$ws.Cells.Item(113, 9).Value = 526.5714
$ws.Cells.Item(113, 10).Value = 631.2727
$ws.Cells.Item(113, 11).Value = 1579.7142
$ws.Cells.Item(113, 12).Value = 1893.8181
$ws.Cells.Item(113, 13).Value = 590.2857999999999
$ws.Cells.Item(113, 14).Value = -6233.8181
$ws.Cells.Item(131, 8).Value = 1283762.9
$ws.Cells.Item(131, 10).Value = 1410159.9
$ws.Cells.Item(131, 12).Value = 4230479.699999999
$ws.Cells.Item(131, 14).Value = -4240559.699999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(14, 8).Value = 4124.9443
$ws.Cells.Item(14, 9).Value = 4124.9443
$ws.Cells.Item(14, 11).Value = 4124.9443
$ws.Cells.Item(14, 13).Value = -3956.9443
$ws.Cells.Item(69, 8).Value = 54000
$ws.Cells.Item(69, 10).Value = 54000
$ws.Cells.Item(69, 12).Value = 54000
$ws.Cells.Item(69, 14).Value = -55498
$ws.Cells.Item(72, 8).Value = 54000
$ws.Cells.Item(72, 10).Value = 54000
$ws.Cells.Item(72, 12).Value = 162000
$ws.Cells.Item(72, 14).Value = -169488

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(20, 8).Value = 3000
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 12).Value = 0
$ws.Cells.Item(20, 14).Value = ""
$ws.Cells.Item(32, 8).Value = 0
$ws.Cells.Item(32, 9).Value = 0
$ws.Cells.Item(32, 11).Value = 0
$ws.Cells.Item(32, 13).Value = ""
$ws.Cells.Item(127, 8).Value = 47821
$ws.Cells.Item(127, 10).Value = 47821
$ws.Cells.Item(127, 12).Value = 47821
$ws.Cells.Item(127, 14).Value = -57741
$ws.Cells.Item(132, 8).Value = 6806.8438
$ws.Cells.Item(132, 9).Value = 7634.7393
$ws.Cells.Item(132, 10).Value = 4691.1113
$ws.Cells.Item(132, 11).Value = 22904.2179
$ws.Cells.Item(132, 12).Value = 14073.3339
$ws.Cells.Item(132, 13).Value = -20374.2179
$ws.Cells.Item(132, 14).Value = -19133.3339

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(17, 8).Value = 12501001
$ws.Cells.Item(17, 9).Value = 12501001
$ws.Cells.Item(17, 11).Value = 12501001
$ws.Cells.Item(17, 13).Value = -12500829
$ws.Cells.Item(23, 8).Value = 6833.3335
$ws.Cells.Item(23, 9).Value = 2500
$ws.Cells.Item(23, 10).Value = 9000
$ws.Cells.Item(23, 11).Value = 2500
$ws.Cells.Item(23, 12).Value = 9000
$ws.Cells.Item(23, 13).Value = -2271
$ws.Cells.Item(23, 14).Value = -9458
$ws.Cells.Item(24, 8).Value = 5000
$ws.Cells.Item(24, 10).Value = 5000
$ws.Cells.Item(24, 12).Value = 5000
$ws.Cells.Item(24, 14).Value = -5460
$ws.Cells.Item(33, 8).Value = 9604
$ws.Cells.Item(33, 10).Value = 9604
$ws.Cells.Item(33, 12).Value = 9604
$ws.Cells.Item(33, 14).Value = -10104
$ws.Cells.Item(36, 8).Value = 9604
$ws.Cells.Item(36, 10).Value = 9604
$ws.Cells.Item(36, 12).Value = 9604
$ws.Cells.Item(36, 14).Value = -10104

